$d = $word.ActiveDocument

# Locate the paragraph: "Another addition to section 2. 2nd Commit."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Another addition to section 2.*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Compute character offsets (relative to the whole story) of the three
# pieces of text inside the paragraph, based on the paragraph's known,
# literal starting text:
#   "Another addition to section 2. " + "2" + "nd" + " Commit."
$prefix = "Another addition to section 2. "
$digit  = "2"
$ord    = "nd"

$prefixStart = $pStart
$digitStart  = $prefixStart + $prefix.Length
$ordStart    = $digitStart + $digit.Length
$ordEnd      = $ordStart + $ord.Length

# 1) Drop the trailing "2" that used to precede the superscripted "nd"
#    (this is a pure in-run text edit, so it does not disturb the
#    neighbouring runs).
$rDigit = $d.Range($digitStart, $digitStart + $digit.Length)
$rDigit.Text = ""

# 2) The old superscript run ("nd") now starts where the digit used to
#    start (it shifted left by one character). Replace its text with
#    "Second" - this keeps it as its own run (only the text inside the
#    run is being edited).
$ordStart2 = $digitStart
$ordEnd2   = $ordStart2 + $ord.Length
$rOrd = $d.Range($ordStart2, $ordEnd2)
$rOrd.Text = "Second"

# 3) Strip the superscript formatting from that run now that it reads
#    "Second".
$newWordStart = $ordStart2
$newWordEnd   = $newWordStart + "Second".Length
$rWord = $d.Range($newWordStart, $newWordEnd)
$rWord.Font.Superscript = 0

# 4) Move the (hidden) "_GoBack" bookmark so that it sits between
#    "Second" and " Commit." instead of after " Commit.".
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    $d.Bookmarks.Add("_GoBack", $d.Range($newWordEnd, $newWordEnd))
}
